# Apply updated crypto price/volume data per commit "Updated cryptos list on Mon Oct 30 08:47:54 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "'34.386.57"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "'  +0.53%  "
$ws.Cells.Item(2,5).Style = "Normal"

# Row 3
$ws.Cells.Item(3,4).Value = "'1.815.17"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "'  +1.49%  "
$ws.Cells.Item(3,5).Style = "Normal"

# Row 4
$ws.Cells.Item(4,5).Value = "'  -0.13%  "
$ws.Cells.Item(4,5).Style = "Normal"

# Row 5
$ws.Cells.Item(5,4).Value = "'228.01"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "'  +0.85%  "
$ws.Cells.Item(5,5).Style = "Normal"

# Row 6
$ws.Cells.Item(6,4).Value = "'0.557"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "'  +1.51%  "
$ws.Cells.Item(6,5).Style = "Normal"

# Row 7
$ws.Cells.Item(7,4).Value = "'0.998"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = "'  -0.27%  "
$ws.Cells.Item(7,5).Style = "Normal"

# Row 8
$ws.Cells.Item(8,4).Value = "'33.80"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "'  +4.43%  "
$ws.Cells.Item(8,5).Style = "Normal"

# Row 9
$ws.Cells.Item(9,4).Value = "'0.298"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "'  +0.93%  "
$ws.Cells.Item(9,5).Style = "Normal"

# Row 10
$ws.Cells.Item(10,4).Value = "'0.0692"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "'  +0.30%  "
$ws.Cells.Item(10,5).Style = "Normal"

# Row 11
$ws.Cells.Item(11,4).Value = "'0.0948"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "'  +0.00%  "
$ws.Cells.Item(11,5).Style = "Normal"

# Row 12
$ws.Cells.Item(12,4).Value = "'2.066.59"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "'  +0.98%  "
$ws.Cells.Item(12,5).Style = "Normal"

# Row 13
$ws.Cells.Item(13,2).Value = "'Chainlink"
$ws.Cells.Item(13,2).Style = "Normal"
$ws.Cells.Item(13,3).Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(13,3).Style = "Normal"
$ws.Cells.Item(13,4).Value = "'11.24"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "'  +1.17%  "
$ws.Cells.Item(13,5).Style = "Normal"

# Row 14
$ws.Cells.Item(14,2).Value = "'WrappedEther"
$ws.Cells.Item(14,2).Style = "Normal"
$ws.Cells.Item(14,3).Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14,3).Style = "Normal"
$ws.Cells.Item(14,4).Value = "'1.805.25"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "'  +0.74%  "
$ws.Cells.Item(14,5).Style = "Normal"

# Row 15
$ws.Cells.Item(15,4).Value = "'0.640"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "'  +2.45%  "
$ws.Cells.Item(15,5).Style = "Normal"

# Row 16
$ws.Cells.Item(16,4).Value = "'34.399.91"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = "'  +0.63%  "
$ws.Cells.Item(16,5).Style = "Normal"

# Row 17
$ws.Cells.Item(17,4).Value = "'4.30"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = "'  +2.28%  "
$ws.Cells.Item(17,5).Style = "Normal"

# Row 18
$ws.Cells.Item(18,4).Value = "'68.71"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "'  +1.08%  "
$ws.Cells.Item(18,5).Style = "Normal"

# Row 19
$ws.Cells.Item(19,4).Value = "'0.0₃0798"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "'  -1.36%  "
$ws.Cells.Item(19,5).Style = "Normal"

# Row 20
$ws.Cells.Item(20,4).Value = "'245.46"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "'  -0.29%  "
$ws.Cells.Item(20,5).Style = "Normal"

# Row 21
$ws.Cells.Item(21,4).Value = "'11.38"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "'  +3.35%  "
$ws.Cells.Item(21,5).Style = "Normal"

# Row 22
$ws.Cells.Item(22,5).Value = "'  +0.12%  "
$ws.Cells.Item(22,5).Style = "Normal"

# Row 23
$ws.Cells.Item(23,4).Value = "'4.18"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "'  +0.31%  "
$ws.Cells.Item(23,5).Style = "Normal"

# Row 24
$ws.Cells.Item(24,4).Value = "'167.52"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = "'  +3.45%  "
$ws.Cells.Item(24,5).Style = "Normal"

# Row 25
$ws.Cells.Item(25,4).Value = "'2.08"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "'  +1.17%  "
$ws.Cells.Item(25,5).Style = "Normal"

# Row 26
$ws.Cells.Item(26,4).Value = "'7.36"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = "'  +2.47%  "
$ws.Cells.Item(26,5).Style = "Normal"

# Row 27
$ws.Cells.Item(27,4).Value = "'16.77"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "'  +2.65%  "
$ws.Cells.Item(27,5).Style = "Normal"

# Row 28
$ws.Cells.Item(28,5).Value = "'  +0.99%  "
$ws.Cells.Item(28,5).Style = "Normal"

# Row 29
$ws.Cells.Item(29,5).Value = "'  -0.42%  "
$ws.Cells.Item(29,5).Style = "Normal"

# Row 30
$ws.Cells.Item(30,4).Value = "'4.00"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = "'  +6.67%  "
$ws.Cells.Item(30,5).Style = "Normal"

# Row 31
$ws.Cells.Item(31,4).Value = "'0.0529"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = "'  +1.45%  "
$ws.Cells.Item(31,5).Style = "Normal"

# Row 32
$ws.Cells.Item(32,5).Value = "'  +1.02%  "
$ws.Cells.Item(32,5).Style = "Normal"

# Row 33
$ws.Cells.Item(33,4).Value = "'3.82"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = "'  +1.59%  "
$ws.Cells.Item(33,5).Style = "Normal"

# Row 34
$ws.Cells.Item(34,4).Value = "'1.84"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "'  +1.73%  "
$ws.Cells.Item(34,5).Style = "Normal"

# Row 35
$ws.Cells.Item(35,5).Value = "'  +2.24%  "
$ws.Cells.Item(35,5).Style = "Normal"

# Row 36
$ws.Cells.Item(36,4).Value = "'1.412.31"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = "'  -2.08%  "
$ws.Cells.Item(36,5).Style = "Normal"

# Row 37
$ws.Cells.Item(37,4).Value = "'0.677"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = "'  +1.84%  "
$ws.Cells.Item(37,5).Style = "Normal"

# Row 38
$ws.Cells.Item(38,5).Value = "'  +1.04%  "
$ws.Cells.Item(38,5).Style = "Normal"

# Row 39
$ws.Cells.Item(39,4).Value = "'0.0191"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "'  +0.00%  "
$ws.Cells.Item(39,5).Style = "Normal"

# Row 40
$ws.Cells.Item(40,4).Value = "'85.77"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "'  +4.26%  "
$ws.Cells.Item(40,5).Style = "Normal"

# Row 41
$ws.Cells.Item(41,4).Value = "'2.84"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "'  +4.77%  "
$ws.Cells.Item(41,5).Style = "Normal"

# Row 42
$ws.Cells.Item(42,2).Value = "'ARBITRUM"
$ws.Cells.Item(42,2).Style = "Normal"
$ws.Cells.Item(42,3).Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(42,3).Style = "Normal"
$ws.Cells.Item(42,4).Value = "'0.954"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "'  +3.39%  "
$ws.Cells.Item(42,5).Style = "Normal"

# Row 43
$ws.Cells.Item(43,2).Value = "'HuobiToken"
$ws.Cells.Item(43,2).Style = "Normal"
$ws.Cells.Item(43,3).Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(43,3).Style = "Normal"
$ws.Cells.Item(43,4).Value = "'2.41"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "'  +0.94%  "
$ws.Cells.Item(43,5).Style = "Normal"

# Row 44
$ws.Cells.Item(44,4).Value = "'14.13"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "'  +1.79%  "
$ws.Cells.Item(44,5).Style = "Normal"

# Row 45
$ws.Cells.Item(45,5).Value = "'  +0.53%  "
$ws.Cells.Item(45,5).Style = "Normal"

# Row 46
$ws.Cells.Item(46,5).Value = "'  +3.01%  "
$ws.Cells.Item(46,5).Style = "Normal"

# Row 47
$ws.Cells.Item(47,4).Value = "'6.06"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "'  -0.60%  "
$ws.Cells.Item(47,5).Style = "Normal"

# Row 48
$ws.Cells.Item(48,4).Value = "'1.965.21"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "'  +0.97%  "
$ws.Cells.Item(48,5).Style = "Normal"

# Row 49
$ws.Cells.Item(49,4).Value = "'105.66"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "'  +0.19%  "
$ws.Cells.Item(49,5).Style = "Normal"

# Row 50
$ws.Cells.Item(50,4).Value = "'0.999"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "'  -0.18%  "
$ws.Cells.Item(50,5).Style = "Normal"

# Row 51
$ws.Cells.Item(51,4).Value = "'0.0₆0126"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "'  -2.29%  "
$ws.Cells.Item(51,5).Style = "Normal"

